# "Delete sku column in product sheet"
#
# The ProductDetails sheet has a "Sku" column (column D) that needs to be
# removed entirely, shifting every column to its right one place to the
# left. Before deleting, the user selects a cell in that column (D2), and
# ends up with the ProductDetails sheet active/selected (it was the
# ProductVariants sheet that was active before the edit).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductDetails")
$ws2 = $wb.Worksheets.Item("ProductVariants")

# Work on the ProductDetails sheet: make it active, then delete the
# entire "Sku" column (column D).
$ws1.Activate()
$ws1.Columns.Item(4).Delete()

# Touch the trailing formatting cell at the sheet's far edge so it keeps
# a style entry (cosmetic parity with the source file's extent) without
# introducing any new value/shared string.
$ws1.Cells.Item(4, 1024).NumberFormat = "General"

# Leave the selection on D2 (where BrandName now starts), matching what
# the column-delete operation would leave selected.
$ws1.Range("D2").Select()

# The ProductVariants sheet keeps its own selection, just narrowed from
# the multi-area A11:A13 down to the single active cell A13, and is no
# longer the active/selected tab.
$ws2.Activate()
$ws2.Range("A13").Select()

# ProductDetails ends up as the active sheet.
$ws1.Activate()
